$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (PEARSON row): "< -0.0745" -> "< -0.0744"
$ws.Range("D2:H2").Value = "< -0.0744"

# Row 2 (PEARSON row): "0.0763 -0.0745" -> "0.077 -0.0744"
$ws.Range("J2:O2").Value = "0.077 -0.0744"

# Row 3 (FISHER row): "0.0931 -0.0463" -> "0.0941 -0.0463"
$ws.Range("J3:O3").Value = "0.0941 -0.0463"

# I4:I8 column: "0.9957 0.0101" -> "0.9958 0.01"
$ws.Range("I4:I8").Value = "0.9958 0.01"

# J4:O8 block: "0.3641 0" -> "0.3677 0"
$ws.Range("J4:O8").Value = "0.3677 0"

# Row 9: "0.2259 -0.0101" -> "0.2281 -0.01"
$ws.Range("J9:O9").Value = "0.2281 -0.01"
